$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("RegistryCode", "Age", "WSO", "WPM_log", "FPM_log")
    ,@("sub 1419", 66, 823.4285714285714, 6.868175328897202, 7.027466179708268)
    ,@("sub 3058", 61, 1422.428571428571, 7.297207142434806, 7.532548481663499)
    ,@("sub 2917", 59, 450.8571428571428, 6.754401649846619, 7.004378343804022)
    ,@("sub 3104", 50, 405.5714285714286, 6.680493671059009, 7.049660684150259)
    ,@("sub 3889", 62, 202.7142857142857, 7.288239991158792, 7.489944295661041)
    ,@("sub 3642", 60, 275, 6.952684029225537, 7.172272124877588)
    ,@("sub 3035", 67, 843.2857142857143, 7.145163944947113, 7.401780564069782)
    ,@("sub 4182", 53, 146.2857142857143, 7.105047066617693, 7.34600224492737)
    ,@("sub 4137", 70, 162.2857142857143, 6.724701539314093, 6.876484644744835)
    ,@("sub 3583", 51, 481.5714285714286, 7.014093534318485, 7.061145858389295)
    ,@("sub 4281", 46, 128.5714285714286, 6.906132341109006, 7.110800607316067)
    ,@("sub 4191", 55, 131.4285714285714, 6.389721196913212, 6.684503026107953)
    ,@("sub 3201", 51, 436.4285714285714, 6.976161311646197, 7.153210769247173)
    ,@("sub 4275", 55, 121.5714285714286, 6.529226008764224, 6.787137721386846)
    ,@("sub 4170", 58, 308, 6.952293994294211, 7.223429300180174)
    ,@("sub 3676", 63, 280.8571428571428, 6.897691652237548, 7.180615986678988)
    ,@("sub 4431", 71, 96.71428571428571, 6.969570368381352, 7.119417519395422)
    ,@("sub 4294", 67, 107.5714285714286, 7.235921224281354, 7.336714164120068)
    ,@("sub 3650", 62, 331.5714285714286, 7.042858499018744, 7.365863212928386)
    ,@("sub 4140", 44, 161.7142857142857, 6.879172072121193, 7.171148892699564)
    ,@("sub 3912", 55, 229.8571428571429, 6.992008888570469, 7.226634189344773)
    ,@("sub 4051", 59, 196.7142857142857, 7.065195922250192, 7.267587259548748)
    ,@("sub 3286", 76, 390.1428571428572, 6.862102357192216, 7.156536197605726)
    ,@("sub 4208", 51, 137.2857142857143, 6.672539146433658, 7.027823283806345)
    ,@("sub 1536", 69, 1061, 6.77639296078531, 7.010622525090596)
    ,@("sub 4439", 56, 79.14285714285714, 7.14185139731475, 7.370452568701144)
    ,@("sub 4466", 65, 95.85714285714286, 6.654818846307865, 6.900469258295833)
    ,@("sub 4504", 55, 78.14285714285714, 7.134752811483845, 7.223234384019707)
    ,@("sub 3887", 68, 310.4285714285714, 6.805950339526814, 6.970406462786104)
    ,@("sub 4391", 61, 101, 6.823411103516777, 7.111104202968168)
    ,@("sub 4546", 40, 123.5714285714286, 6.574279595328631, 6.811365395628536)
    ,@("sub 4227", 38, 186.7142857142857, 6.421239602538288, 6.607473841627368)
    ,@("sub 4532", 63, 89.42857142857143, 6.961060523195856, 7.138656893841168)
    ,@("sub 3299", 60, 509.2857142857143, 6.938973119086882, 7.190838557880011)
    ,@("sub 4720", 59, 43.28571428571428, 6.959349987793811, 7.083148125701578)
    ,@("sub 4607", 63, 60.42857142857143, 7.225240799600637, 7.358358766887639)
    ,@("sub 4728", 43, 34.42857142857143, 6.858551904453062, 7.137558618980075)
    ,@("sub 4769", 54, 31.71428571428572, 7.342688249360188, 7.50900312285984)
    ,@("sub 4777", 51, 28.71428571428572, 6.852448473005276, 7.045906303804517)
    ,@("sub 3396", 54, 337.2857142857143, 6.765143826059799, 7.05851477048068)
    ,@("sub 4467", 60, 31.28571428571428, 7.07880722235384, 7.232984093428745)
    ,@("sub 3392", 53, 331.1428571428572, 6.941884155244977, 7.179055921613601)
    ,@("sub 3960", 79, 239.7142857142857, 7.330195491296263, 7.282435106535899)
    ,@("sub 3154", 30, 443.1428571428572, 6.949298581924364, 7.143466335843881)
    ,@("sub 3994", 55, 179.8571428571429, 6.941965574200417, 7.265568036903425)
    ,@("sub 4086", 59, 141.7142857142857, 6.983195060598846, 7.296772856863443)
    ,@("sub 4002", 69, 185.2857142857143, 7.224503660714746, 7.374136614896833)
    ,@("sub 2522", 40, 555.8571428571429, 6.711985204335674, 7.05124192015805)
    ,@("sub 3364", 34, 334.4285714285714, 6.964977817750166, 7.158015818784285)
    ,@("sub 2884", 55, 516.8571428571429, 6.847440856796569, 7.138000824057674)
    ,@("sub 2998", 59, 469.2857142857143, 7.282931011531638, 7.576867577619332)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}
